# "remove all hardcoded test data" -
#   - offlineShopingProcess: insert a new "state" column (with value
#     "Alabama") before the existing zipCode/phone columns.
#   - add a new "addToCartAndVerifyPrice" sheet (between
#     offlineShopingProcess and signIn) holding size/color/qnt test data.

$wb = $excel.ActiveWorkbook

# --- offlineShopingProcess: insert column G ("state" / "Alabama") -----
$ws1 = $wb.Worksheets.Item("offlineShopingProcess")
$ws1.Columns("G:G").Insert()
$ws1.Range("G1").Value = "state"
$ws1.Range("G2").Value = "Alabama"

# --- new sheet: addToCartAndVerifyPrice --------------------------------
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "addToCartAndVerifyPrice"
$ws2.Range("A1").Value = "size"
$ws2.Range("B1").Value = "color"
$ws2.Range("C1").Value = "qnt"
$ws2.Range("B2").Value = "Blue"
$ws2.Range("A2").Value = "L"
$ws2.Range("C2").Value = 6

# --- restore per-sheet selections / active sheet -----------------------
[void]$ws1.Activate()
[void]$ws1.Range("L6").Select()

[void]$ws2.Activate()
[void]$ws2.Range("A2").Select()

# --- best-effort: reposition the workbook window (cosmetic) ------------
try {
    $excel.ActiveWindow.Left = 3555
    $excel.ActiveWindow.Top = 5280
} catch {
}
